$wb = $excel.ActiveWorkbook

# --- Sheet: DatosCuenta ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokQAJuneLastOne"
$wsCuenta.Range("B2").Value = "SmokeNameQAJuneLastOne"
$wsCuenta.Range("C2").Value = 27100132
$wsCuenta.Range("D2").Value = 133

# --- Sheet: DatosHogar ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 652

# --- Sheet: DatosMotor ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP034"
$wsMotor.Range("B2").Value = "ABC12SSMP034"
$wsMotor.Range("C2").Value = "ZAZ123SSMP034"

# --- Sheet: DatosAP ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200133
$wsAP.Activate()
$wsAP.Range("A3").Select()
